$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.294.69"
$ws.Range("E2").Value = "  -7.06%  "
$ws.Range("D3").Value = "3.542.03"
$ws.Range("E3").Value = "  -3.22%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "390.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -7.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.46"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.49%  "
$ws.Range("D7").Value = "3.533.63"
$ws.Range("E7").Value = "  -3.07%  "
$ws.Range("E8").Value = "  -10.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  -11.95%  "
$ws.Range("E11").Value = "  -23.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000322"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -27.43%  "
$ws.Range("E13").Value = "  -8.22%  "
$ws.Range("D14").Value = "4.093.23"
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.11"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -7.13%  "
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.24"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.13%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.523.98"
$ws.Range("E18").Value = "  -3.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -7.07%  "
$ws.Range("D20").Value = "63.377.75"
$ws.Range("E20").Value = "  -6.75%  "
$ws.Range("E21").Value = "  -9.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.89"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -14.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.84"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.93%  "
$ws.Range("E25").Value = "  -4.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +9.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.54"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.97"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.64"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -14.10%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.65"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.72"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.110"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.72"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.00%  "
$ws.Range("E34").Value = "  -6.13%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.38"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -9.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.28"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0434"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -11.82%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "0.0₃0638"
$ws.Range("E40").Value = "  -19.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.66"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.08"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +16.25%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.130"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -13.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "139.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.17%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.31"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +18.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.46"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -9.07%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.06"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.64"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -9.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.273"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -9.38%  "
